$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 14; this shifts existing rows 14-61 down to 15-62
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the new record's data
$ws.Cells.Item(14, 1).Value = 5
$ws.Cells.Item(14, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(14, 3).Value = "Maule"
$ws.Cells.Item(14, 4).Value = 44481
$ws.Cells.Item(14, 5).Value = 7
$ws.Cells.Item(14, 6).Value = 100112013
$ws.Cells.Item(14, 7).Value = "Alcachofa"
$ws.Cells.Item(14, 8).Value = "Madrigal"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 300
$ws.Cells.Item(14, 11).Value = 10000
$ws.Cells.Item(14, 12).Value = 10000
$ws.Cells.Item(14, 13).Value = 10000
$ws.Cells.Item(14, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(14, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(14, 16).Value = 250
$ws.Cells.Item(14, 17).Value = 40
$ws.Cells.Item(14, 18).Value = "Hortaliza"

# Ensure date formatting/style for column D in the new row matches the rest of the column
$ws.Cells.Item(14, 4).NumberFormat = $ws.Cells.Item(15, 4).NumberFormat
